# subj3/results_Alpha.xlsx — combine movement-condition rows so a
# two-sample t-test (ttest_ind) can be run on pooled "eyescrunching+jaw"
# vs "jaw+raisingeyebrows" groups instead of the five separate conditions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded condition rows (movehat / movehead and the
# separate raisingeyebrows row), leaving just the header plus the two
# pooled-condition rows.
$ws.Rows("4:6").Delete()

# Relabel the two remaining condition rows to reflect the pooled groups.
$ws.Range("U2").Value = "eyescrunching+jaw"
$ws.Range("U3").Value = "jaw+raisingeyebrows"

# Updated SNR LMS stats (Q/R/S) recomputed for the pooled groups.
$ws.Range("Q2").Value = -0.5382516940701341
$ws.Range("R2").Value = 9.111637823448699
$ws.Range("S2").Value = -11.93809941404094

$ws.Range("Q3").Value = -0.3122190014324168
$ws.Range("R3").Value = 11.76986702342861
$ws.Range("S3").Value = 4.459231447394495
